# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet: bump the timestamp, update the
# case/recovered/death counters for several countries, and fix a handful
# of countries that were associated with the wrong row (their counters
# follow the name to its corrected row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 17:22"

$ws.Range("B4").Value = 321337
$ws.Range("C4").Value = 9980
$ws.Range("D4").Value = 16553
$ws.Range("E4").Value = 295656
$ws.Range("F4").Value = 8468
$ws.Range("G4").Value = 676
$ws.Range("H4").Value = 9128

$ws.Range("B7").Value = 97351
$ws.Range("C7").Value = 1259
$ws.Range("E7").Value = 69472
$ws.Range("G7").Value = 35
$ws.Range("H7").Value = 1479

$ws.Range("B17").Value = 11930
$ws.Range("C17").Value = 149
$ws.Range("E17").Value = 8728

$ws.Range("B19").Value = 10475
$ws.Range("C19").Value = 115
$ws.Range("E19").Value = 9900
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 448

$ws.Range("E23").Value = 3337
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 35

$ws.Range("A66").Value = "Moldavia"
$ws.Range("B66").Value = 864
$ws.Range("C66").Value = 112
$ws.Range("D66").Value = 30
$ws.Range("E66").Value = 819
$ws.Range("F66").Value = 80
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 15

$ws.Range("A67").Value = "Armenia"
$ws.Range("B67").Value = 822
$ws.Range("C67").Value = 52
$ws.Range("D67").Value = 57
$ws.Range("E67").Value = 758
$ws.Range("F67").Value = 30
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 7

$ws.Range("A68").Value = "Lituania"
$ws.Range("B68").Value = 811
$ws.Range("C68").Value = 40
$ws.Range("D68").Value = 7
$ws.Range("E68").Value = 792
$ws.Range("F68").Value = 11
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 12

$ws.Range("A85").Value = "Republica de Chipre"
$ws.Range("B85").Value = 446
$ws.Range("C85").Value = 20
$ws.Range("D85").Value = 33
$ws.Range("E85").Value = 404
$ws.Range("F85").Value = 11
$ws.Range("H85").Value = 9

$ws.Range("A86").Value = "Costa Rica"
$ws.Range("B86").Value = 435
$ws.Range("D86").Value = 13
$ws.Range("E86").Value = 420
$ws.Range("F86").Value = 13
$ws.Range("H86").Value = 2

$ws.Range("A156").Value = "Birmania"
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 0
$ws.Range("H156").Value = 1

$ws.Range("A157").Value = "Haiti"
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 1
$ws.Range("H157").Value = 0

$ws.Range("A176").Value = "Sudan"
$ws.Range("B176").Value = 12
$ws.Range("C176").Value = 2
$ws.Range("D176").Value = 2
$ws.Range("E176").Value = 8
$ws.Range("H176").Value = 2

$ws.Range("A177").Value = "Laos"
$ws.Range("C177").Value = 1
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 11

$ws.Range("A178").Value = "Groenlandia"
$ws.Range("D178").Value = 3
$ws.Range("E178").Value = 8
$ws.Range("H178").Value = 0

$ws.Range("A179").Value = "Curazao"
$ws.Range("B179").Value = 11
$ws.Range("D179").Value = 5
$ws.Range("E179").Value = 5
$ws.Range("H179").Value = 1

$ws.Range("A180").Value = "San Cristobal y Nieves"
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 10

$ws.Range("A181").Value = "Seychelles"
$ws.Range("E181").Value = 10
$ws.Range("H181").Value = 0

$ws.Range("A182").Value = "Surinam"
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 9
$ws.Range("H182").Value = 1

$ws.Range("A183").Value = "Mozambique"
$ws.Range("D183").Value = 1
$ws.Range("E183").Value = 9
$ws.Range("H183").Value = 0

$ws.Range("A184").Value = "Angola"
$ws.Range("B184").Value = 10
$ws.Range("D184").Value = 2
$ws.Range("E184").Value = 6
$ws.Range("H184").Value = 2

$ws.Range("A187").Value = "Zimbabue"
$ws.Range("D187").Value = 0
$ws.Range("H187").Value = 1

$ws.Range("A188").Value = "Nepal"
$ws.Range("D188").Value = 1
$ws.Range("H188").Value = 0

$ws.Range("A192").Value = "Cabo Verde"
$ws.Range("D192").Value = 0
$ws.Range("H192").Value = 1

$ws.Range("A194").Value = "San Vicente y las Granadinas"
$ws.Range("D194").Value = 1
$ws.Range("H194").Value = 0

$ws.Range("A198").Value = "Belice"
$ws.Range("C198").Value = 1
$ws.Range("F198").Value = 1

$ws.Range("A199").Value = "Nicaragua"
$ws.Range("C199").Value = 0
$ws.Range("E199").Value = 4
$ws.Range("F199").Value = 0
$ws.Range("H199").Value = 1

$ws.Range("A200").Value = "Islas Turcas y Caicos"
$ws.Range("G200").Value = 1

$ws.Range("A202").Value = "Malaui"

$ws.Range("A203").Value = "Sahara Occidental"
